$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2796
$ws.Range("L3").Value = 2839
$ws.Range("K4").Value = 1767
$ws.Range("L4").Value = 760
$ws.Range("L6").Value = 2541
$ws.Range("K7").Value = 27559
$ws.Range("L7").Value = 9099

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 306
$ws.Range("L8").Value = 575
$ws.Range("L13").Value = 14
$ws.Range("L15").Value = 66
$ws.Range("L19").Value = 255
$ws.Range("L20").Value = 233
$ws.Range("L21").Value = 25
$ws.Range("L23").Value = 92
$ws.Range("L25").Value = 51
$ws.Range("L29").Value = 491
$ws.Range("L31").Value = 87
$ws.Range("L33").Value = 414
$ws.Range("L37").Value = 337
$ws.Range("L38").Value = 8
$ws.Range("L41").Value = 42
$ws.Range("L43").Value = 70
$ws.Range("L46").Value = 20
$ws.Range("L48").Value = 122
$ws.Range("L49").Value = 49
$ws.Range("L50").Value = 49
$ws.Range("L51").Value = 111
$ws.Range("L54").Value = 185
$ws.Range("L59").Value = 12
$ws.Range("L61").Value = 12
$ws.Range("K63").Value = 159
$ws.Range("L63").Value = 28
$ws.Range("L64").Value = 57
$ws.Range("L67").Value = 336
$ws.Range("L69").Value = 26
$ws.Range("L79").Value = 244
$ws.Range("L85").Value = 463
$ws.Range("L89").Value = 117
$ws.Range("L90").Value = 90
$ws.Range("L91").Value = 129
$ws.Range("L92").Value = 27
$ws.Range("L93").Value = 45
$ws.Range("L99").Value = 152
$ws.Range("K101").Value = 27559
$ws.Range("L101").Value = 9099

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L4").Value = 26
$ws.Range("L6").Value = 84
$ws.Range("L7").Value = 306

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 135
$ws.Range("L3").Value = 188
$ws.Range("L7").Value = 463

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 162
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 575

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 115
$ws.Range("L6").Value = 143
$ws.Range("L7").Value = 414

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 98
$ws.Range("L3").Value = 98
$ws.Range("L6").Value = 108
$ws.Range("L7").Value = 337

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 37
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 100
$ws.Range("L7").Value = 336

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 37
$ws.Range("L6").Value = 93
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 182
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 491

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 86
$ws.Range("L3").Value = 80
$ws.Range("L7").Value = 255

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 7
$ws.Range("L6").Value = 14

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 81
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 71
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 32
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 8
